$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 12, pushing existing rows 12:33 down to 13:34.
$ws.Rows.Item(12).Insert()

# Populate the newly inserted row 12 with the new Cilantro price record.
$ws.Range("A12").Value = 5
$ws.Range("B12").Value = "Macroferia Regional de Talca"
$ws.Range("C12").Value = "Maule"
$ws.Range("D12").Value = 44757
$ws.Range("E12").Value = 7
$ws.Range("F12").Value = 100112040
$ws.Range("G12").Value = "Cilantro"
$ws.Range("H12").Value = "Sin especificar"
$ws.Range("I12").Value = "Primera"
$ws.Range("J12").Value = 200
$ws.Range("K12").Value = 14000
$ws.Range("L12").Value = 14000
$ws.Range("M12").Value = 14000
$ws.Range("N12").Value = "$/caja 36 atados"
$ws.Range("O12").Value = "Región Metropolitana"
$ws.Range("P12").Value = 389
$ws.Range("Q12").Value = 36
$ws.Range("R12").Value = "Hortaliza"
